$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 4 (bug #2): reassign to Mukesh, clear Status and Developer comment
$ws.Range("C4").Value = "Mukesh"
$ws.Range("H4").ClearContents()
$ws.Range("I4").ClearContents()
$ws.Rows("4").AutoFit()

# Row 5 (bug #3): mark Resolved
$ws.Range("H5").Value = "Resolved"

# Fix sheet view / selection
$ws.Range("I5").Select()
